$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 617, shifting the existing data (rows 617:677) down to 619:679
$ws.Rows("617:618").Insert()

# Fill in the new row 617 (copy static columns from the row below, which is the
# shifted original row 617 data now sitting at row 619)
$ws.Cells.Item(617, 1).Value = 9
$ws.Cells.Item(617, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(617, 3).Value = "Metropolitana"
$ws.Cells.Item(617, 4).Value = 44769
$ws.Cells.Item(617, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(617, 5).Value = 13
$ws.Cells.Item(617, 6).Value = 100112040
$ws.Cells.Item(617, 7).Value = "Cilantro"
$ws.Cells.Item(617, 8).Value = "Sin especificar"
$ws.Cells.Item(617, 9).Value = "Primera"
$ws.Cells.Item(617, 10).Value = 70
$ws.Cells.Item(617, 11).Value = 12000
$ws.Cells.Item(617, 12).Value = 12000
$ws.Cells.Item(617, 13).Value = 12000
$ws.Cells.Item(617, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(617, 15).Value = "Región Metropolitana"
$ws.Cells.Item(617, 16).Value = 333
$ws.Cells.Item(617, 17).Value = 36
$ws.Cells.Item(617, 18).Value = "Hortaliza"

# Fill in the new row 618
$ws.Cells.Item(618, 1).Value = 9
$ws.Cells.Item(618, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(618, 3).Value = "Metropolitana"
$ws.Cells.Item(618, 4).Value = 44769
$ws.Cells.Item(618, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(618, 5).Value = 13
$ws.Cells.Item(618, 6).Value = 100112040
$ws.Cells.Item(618, 7).Value = "Cilantro"
$ws.Cells.Item(618, 8).Value = "Sin especificar"
$ws.Cells.Item(618, 9).Value = "Primera"
$ws.Cells.Item(618, 10).Value = 160
$ws.Cells.Item(618, 11).Value = 18000
$ws.Cells.Item(618, 12).Value = 20000
$ws.Cells.Item(618, 13).Value = 19000
$ws.Cells.Item(618, 14).Value = "`$/docena de atados"
$ws.Cells.Item(618, 15).Value = "Región Metropolitana"
$ws.Cells.Item(618, 16).Value = 6333
$ws.Cells.Item(618, 17).Value = 3
$ws.Cells.Item(618, 18).Value = "Hortaliza"
